$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.00"
$ws.Range("E2").Value = "'2.12%"

$ws.Range("D3").Value = "'41.42"
$ws.Range("E3").Value = "'2.53%"

$ws.Range("D4").Value = "'5.034"
$ws.Range("E4").Value = "'-0.42%"

$ws.Range("D5").Value = "'0.07554"
$ws.Range("E5").Value = "'3.73%"

$ws.Range("B6").Value = "'GateToken"
$ws.Range("C6").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.387"
$ws.Range("E6").Value = "'2.37%"

$ws.Range("B7").Value = "'FTXToken"
$ws.Range("C7").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.603"
$ws.Range("E7").Value = "'1.83%"

$ws.Range("B8").Value = "'MXToken"
$ws.Range("C8").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9281"
$ws.Range("E8").Value = "'0.85%"

$ws.Range("B9").Value = "'BTSEToken"
$ws.Range("C9").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.410"
$ws.Range("E9").Value = "'2.21%"

$ws.Range("B10").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1222"
$ws.Range("E10").Value = "'5.51%"

$ws.Range("B11").Value = "'WazirX"
$ws.Range("C11").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1837"
$ws.Range("E11").Value = "'6.47%"

$ws.Range("B12").Value = "'MandalaExchangeToken"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08891"
$ws.Range("E12").Value = "'3.15%"

$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03990"
$ws.Range("E13").Value = "'-4.83%"

$ws.Range("B14").Value = "'BitMartToken"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1052"
$ws.Range("E14").Value = "'-0.06%"

$ws.Range("B15").Value = "'BitForexToken"
$ws.Range("C15").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001282"
$ws.Range("E15").Value = "'0.60%"

$ws.Range("B16").Value = "'TigerCash"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005958"
$ws.Range("E16").Value = "'2.45%"

$ws.Range("B17").Value = "'LEO"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.336"
$ws.Range("E17").Value = "'-1.86%"

$ws.Range("D18").Value = "'0.3321"
$ws.Range("E18").Value = "'1.34%"

$ws.Range("D19").Value = "'7.945"
$ws.Range("E19").Value = "'1.05%"

$ws.Range("D20").Value = "'0.1420"
$ws.Range("E20").Value = "'4.97%"

$ws.Range("E21").Value = "'3.92%"

$ws.Range("D22").Value = "'0.04064"
$ws.Range("E22").Value = "'4.95%"

$ws.Range("E23").Value = "'-0.40%"

$ws.Range("D24").Value = "'0.003982"
$ws.Range("E24").Value = "'5.14%"

$ws.Range("D25").Value = "'0.0001229"
$ws.Range("E25").Value = "'-4.05%"

$ws.Range("E26").Value = "'-0.12%"

$ws.Range("D38").Value = "'0.02402"
$ws.Range("E38").Value = "'3.80%"

$ws.Range("D39").Value = "'0.05205"
$ws.Range("E39").Value = "'5.16%"

$ws.Range("D40").Value = "'0.006385"
$ws.Range("E40").Value = "'-1.40%"

$ws.Range("D41").Value = "'0.007759"
$ws.Range("E41").Value = "'0.97%"

$ws.Range("E42").Value = "'4.56%"

$ws.Range("D43").Value = "'0.007489"
$ws.Range("E43").Value = "'1.44%"

$ws.Range("D44").Value = "'0.007837"
$ws.Range("E44").Value = "'11.25%"

$ws.Range("D45").Value = "'0.3215"
$ws.Range("E45").Value = "'10.73%"

$ws.Range("D46").Value = "'0.00006788"
$ws.Range("E46").Value = "'5.68%"

$ws.Range("E47").Value = "'-0.05%"

$ws.Range("D48").Value = "'0.04628"
$ws.Range("E48").Value = "'124.29%"

$ws.Range("D49").Value = "'0.004201"
$ws.Range("E49").Value = "'-0.01%"

$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.05%"

$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.05%"
